$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# atividade 4 - modelagem parametrica
# Row 7 = "BIM 3D - Modelagem Parametrica": Unidade 3 (H) and Unidade 4 (I) concluded -> 10 pts each
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 10

# Reflect the cursor position left by the user after the edit
$ws.Range("I3").Select()
